$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip the leading "root-code-name." / "1:" prefix from the ID / Numeric ID
# columns for rows 3-17 (row 2 itself is "root-code-name" / "1" and stays as-is).
$ws.Range("A3").Value = "common-code-name-1"
$ws.Range("B3").Value = "2"

$ws.Range("A4").Value = "common-code-name-1.common-code-name-2"
$ws.Range("B4").Value = "2:3"

$ws.Range("A5").Value = "common-code-name-1.common-code-name-2.common-code-name-3-1"
$ws.Range("B5").Value = "2:3:4"

$ws.Range("A6").Value = "common-code-name-1.common-code-name-2.common-code-name-3-2"
$ws.Range("B6").Value = "2:3:5"

$ws.Range("A7").Value = "common-code-name-1.common-code-name-2.common-code-name-3-3"
$ws.Range("B7").Value = "2:3:6"

$ws.Range("A8").Value = "error-code-name-1"
$ws.Range("B8").Value = "7"

$ws.Range("A9").Value = "error-code-name-1.error-code-name-2"
$ws.Range("B9").Value = "7:8"

$ws.Range("A10").Value = "error-code-name-1.error-code-name-2.error-code-name-3-1"
$ws.Range("B10").Value = "7:8:9"

$ws.Range("A11").Value = "error-code-name-1.error-code-name-2.error-code-name-3-2"
$ws.Range("B11").Value = "7:8:10"

$ws.Range("A12").Value = "error-code-name-1.error-code-name-2.error-code-name-3-3"
$ws.Range("B12").Value = "7:8:11"

$ws.Range("A13").Value = "error-2-code-name-1"
$ws.Range("B13").Value = "12"

$ws.Range("A14").Value = "error-2-code-name-1.error-2-code-name-2"
$ws.Range("B14").Value = "12:13"

$ws.Range("A15").Value = "error-2-code-name-1.error-2-code-name-2.error-2-code-name-3-1"
$ws.Range("B15").Value = "12:13:14"

$ws.Range("A16").Value = "error-2-code-name-1.error-2-code-name-2.error-2-code-name-3-2"
$ws.Range("B16").Value = "12:13:15"

$ws.Range("A17").Value = "error-2-code-name-1.error-2-code-name-2.error-2-code-name-3-3"
$ws.Range("B17").Value = "12:13:16"

# Update the active cell / selection from B17 to B3
$ws.Range("B3").Select()
